# The workbook previously stored several "Названия университетов" (university
# names) cells with a stray "null" text prefix (an artifact of the app's old
# JSON-writing code, per the commit message "Added writing to JSON file").
# Strip that leftover "null" prefix from each affected cell in column E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "null"

for ($r = 2; $r -le 5; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().StartsWith($prefix)) {
        $cell.Value = $val.ToString().Substring($prefix.Length)
    }
}

# The shorter text now needs a narrower, best-fit column width for column E
# (was 139.88671875 characters, now 136.46484375 characters). Set the closest
# reachable width via the ColumnWidth COM property.
$ws.Columns.Item(5).ColumnWidth = 135.66666666666666
